# Add "hydrogen combined cycle" as a new power plant type on the
# "CPPbES" (CPP CO2 Capture Potential by Electricity Source) sheet.
#
# 1. Rename the existing "hydrogen" entry (row 24) to
#    "hydrogen combustion turbine" and give it the new vertically
#    centered / black font style.
# 2. Add a new row 25 "hydrogen combined cycle" with the same style,
#    and a capture-potential value of 0.
# 3. Leave the selection on the CPPbES sheet pointing at D28 (matching
#    where the author clicked next), while keeping the "About" sheet as
#    the active/selected tab of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CPPbES")

# Row 24: "hydrogen" -> "hydrogen combustion turbine"
$ws.Range("A24").Value = "hydrogen combustion turbine"
$ws.Range("A24").Font.Color = 0
$ws.Range("A24").VerticalAlignment = -4108

# Row 25 (new): "hydrogen combined cycle", using the same formatting as
# row 24 (copy the already-built style rather than re-deriving it, so
# we don't leave any unused/duplicate cell formats behind).
$ws.Range("A24").Copy($ws.Range("A25"))
$ws.Range("A25").Value = "hydrogen combined cycle"
$ws.Range("B25").Value = 0

# Update the on-sheet selection to D28, then restore "About" as the
# active sheet/tab (matches the author's saved workbook state).
$ws.Range("D28").Select()
$wb.Worksheets.Item("About").Select()
